$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns AE (Regional Manager Name) and AF (Regional Manager Contact Name) ---
# Header row (row 1) - bold header style already applied to these (previously empty) cells.
$ws.Range("AE1").Value = "Regional Manager Name"
$ws.Range("AF1").Value = "Regional Manager Contact Name"

# Template-token row (row 2) - matches the centered style used by the rest of row 2.
$ws.Range("AE2").Value = "{vendor:sf_rm_name}"
$ws.Range("AF2").Value = "{vendor:sf_rm_phone}"
$ws.Range("AE2:AF2").HorizontalAlignment = -4108

# --- Column widths: split the old shared AE:AF width into two distinct widths ---
$ws.Columns.Item(31).ColumnWidth = 22
$ws.Columns.Item(32).ColumnWidth = 26.33

# --- Update the saved selection to match the author's final cursor position ---
$ws.Range("AD11").Select() | Out-Null
